# Auto-generated Excel COM-interop script
# Applies the cell-value updates described by the commit diff to the
# "Mateus_Profits" profit-tracking workbook (8 crafting-job sheets).
# All target cells are plain numeric literals (no formulas anywhere in the
# workbook), so each change is a simple value assignment through the
# worksheet.Cells(row, col) COM surface.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 4999.4644   # H40: was 5170.893
$ws.Cells.Item(40, 9).Value = 3903   # I40: was 4003.5715
$ws.Cells.Item(40, 10).Value = 5438.05   # J40: was 5560
$ws.Cells.Item(40, 11).Value = 3903   # K40: was 4003.5715
$ws.Cells.Item(40, 12).Value = 5438.05   # L40: was 5560
$ws.Cells.Item(40, 13).Value = -3728   # M40: was -3828.5715
$ws.Cells.Item(40, 14).Value = -5788.05   # N40: was -5910
$ws.Cells.Item(86, 8).Value = 1000   # H86: was 0
$ws.Cells.Item(86, 9).Value = 1000   # I86: was 0
$ws.Cells.Item(86, 11).Value = 1000   # K86: was 0
$ws.Cells.Item(86, 13).Value = 123   # M86: was None
$ws.Cells.Item(89, 8).Value = 1000   # H89: was 0
$ws.Cells.Item(89, 9).Value = 1000   # I89: was 0
$ws.Cells.Item(89, 11).Value = 5000   # K89: was 0
$ws.Cells.Item(89, 13).Value = 616   # M89: was None
$ws.Cells.Item(98, 8).Value = 2209.9167   # H98: was 2264.4856
$ws.Cells.Item(98, 9).Value = 2209.9167   # I98: was 2264.4856
$ws.Cells.Item(98, 11).Value = 2209.9167   # K98: was 2264.4856
$ws.Cells.Item(98, 13).Value = -711.9167000000002   # M98: was -766.4856
$ws.Cells.Item(122, 8).Value = 2209.9167   # H122: was 2264.4856
$ws.Cells.Item(122, 9).Value = 2209.9167   # I122: was 2264.4856
$ws.Cells.Item(122, 11).Value = 6629.750100000001   # K122: was 6793.4568
$ws.Cells.Item(122, 13).Value = -4179.750100000001   # M122: was -4343.4568
# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4859.5103   # H32: was 4940.146
$ws.Cells.Item(32, 9).Value = 4003.9211   # I32: was 4085.4055
$ws.Cells.Item(32, 11).Value = 4003.9211   # K32: was 4085.4055
$ws.Cells.Item(32, 13).Value = -3716.9211   # M32: was -3798.4055
$ws.Cells.Item(39, 8).Value = 20069   # H39: was 17534.5
$ws.Cells.Item(39, 10).Value = 0   # J39: was 15000
$ws.Cells.Item(39, 12).Value = 0   # L39: was 15000
$ws.Cells.Item(39, 14).Value = $null   # N39: clear (was -16040)
$ws.Cells.Item(74, 8).Value = 3785.6765   # H74: was 4009.875
$ws.Cells.Item(74, 9).Value = 2596.375   # I74: was 2814.3635
$ws.Cells.Item(74, 11).Value = 2596.375   # K74: was 2814.3635
$ws.Cells.Item(74, 13).Value = -1722.375   # M74: was -1940.3635
$ws.Cells.Item(76, 8).Value = 100000   # H76: was 0
$ws.Cells.Item(76, 10).Value = 100000   # J76: was 0
$ws.Cells.Item(76, 12).Value = 100000   # L76: was 0
$ws.Cells.Item(76, 14).Value = -100676   # N76: was None
$ws.Cells.Item(77, 8).Value = 3785.6765   # H77: was 4009.875
$ws.Cells.Item(77, 9).Value = 2596.375   # I77: was 2814.3635
$ws.Cells.Item(77, 11).Value = 12981.875   # K77: was 14071.8175
$ws.Cells.Item(77, 13).Value = -8613.875   # M77: was -9703.817499999999
$ws.Cells.Item(79, 8).Value = 100000   # H79: was 0
$ws.Cells.Item(79, 10).Value = 100000   # J79: was 0
$ws.Cells.Item(79, 12).Value = 100000   # L79: was 0
$ws.Cells.Item(79, 14).Value = -102340   # N79: was None
$ws.Cells.Item(97, 8).Value = 666.4583   # H97: was 685.2273
$ws.Cells.Item(97, 9).Value = 615.1177   # I97: was 635.8
$ws.Cells.Item(97, 11).Value = 615.1177   # K97: was 635.8
$ws.Cells.Item(97, 13).Value = -119.1177   # M97: was -139.8
$ws.Cells.Item(110, 8).Value = 6348.0527   # H110: was 6645.4443
$ws.Cells.Item(110, 9).Value = 3393.1667   # I110: was 3611.182
$ws.Cells.Item(110, 11).Value = 3393.1667   # K110: was 3611.182
$ws.Cells.Item(110, 13).Value = -1348.1667   # M110: was -1566.182
$ws.Cells.Item(122, 8).Value = 1863.4117   # H122: was 1897.2354
$ws.Cells.Item(122, 9).Value = 1592.8   # I122: was 1638.8
$ws.Cells.Item(122, 11).Value = 4778.4   # K122: was 4916.4
$ws.Cells.Item(122, 13).Value = -2328.4   # M122: was -2466.4
# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 3926   # H134: was 4000.524
$ws.Cells.Item(134, 9).Value = 3926   # I134: was 4000.524
$ws.Cells.Item(134, 11).Value = 11778   # K134: was 12001.572
$ws.Cells.Item(134, 13).Value = -9243   # M134: was -9466.572
# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2484.1304   # H16: was 2849.2632
$ws.Cells.Item(16, 9).Value = 1827.6428   # I16: was 2015.5834
$ws.Cells.Item(16, 10).Value = 3505.3333   # J16: was 4278.4287
$ws.Cells.Item(16, 11).Value = 1827.6428   # K16: was 2015.5834
$ws.Cells.Item(16, 12).Value = 3505.3333   # L16: was 4278.4287
$ws.Cells.Item(16, 13).Value = -1540.6428   # M16: was -1728.5834
$ws.Cells.Item(16, 14).Value = -4079.3333   # N16: was -4852.4287
$ws.Cells.Item(58, 8).Value = 8184.9   # H58: was 7112.4165
$ws.Cells.Item(58, 9).Value = 4412.5   # I58: was 3810
$ws.Cells.Item(58, 10).Value = 10699.833   # J58: was 9471.286
$ws.Cells.Item(58, 11).Value = 4412.5   # K58: was 3810
$ws.Cells.Item(58, 12).Value = 10699.833   # L58: was 9471.286
$ws.Cells.Item(58, 13).Value = -4209.5   # M58: was -3607
$ws.Cells.Item(58, 14).Value = -11105.833   # N58: was -9877.286
$ws.Cells.Item(107, 8).Value = 2148.12   # H107: was 2168.6155
$ws.Cells.Item(107, 9).Value = 2358.5833   # I107: was 2373.6667
$ws.Cells.Item(107, 10).Value = 1953.8462   # J107: was 1992.8572
$ws.Cells.Item(107, 11).Value = 2358.5833   # K107: was 2373.6667
$ws.Cells.Item(107, 12).Value = 1953.8462   # L107: was 1992.8572
$ws.Cells.Item(107, 13).Value = -438.5832999999998   # M107: was -453.6667000000002
$ws.Cells.Item(107, 14).Value = -5793.8462   # N107: was -5832.8572
$ws.Cells.Item(113, 8).Value = 2484.1304   # H113: was 2849.2632
$ws.Cells.Item(113, 9).Value = 1827.6428   # I113: was 2015.5834
$ws.Cells.Item(113, 10).Value = 3505.3333   # J113: was 4278.4287
$ws.Cells.Item(113, 11).Value = 1827.6428   # K113: was 2015.5834
$ws.Cells.Item(113, 12).Value = 3505.3333   # L113: was 4278.4287
$ws.Cells.Item(113, 13).Value = 342.3571999999999   # M113: was 154.4166
$ws.Cells.Item(113, 14).Value = -7845.3333   # N113: was -8618.4287
$ws.Cells.Item(122, 8).Value = 3136.7083   # H122: was 2774.9312
$ws.Cells.Item(122, 9).Value = 3215.25   # I122: was 2935.1304
$ws.Cells.Item(122, 10).Value = 2744   # J122: was 2160.8333
$ws.Cells.Item(122, 11).Value = 9645.75   # K122: was 8805.3912
$ws.Cells.Item(122, 12).Value = 8232   # L122: was 6482.499899999999
$ws.Cells.Item(122, 13).Value = -7195.75   # M122: was -6355.3912
$ws.Cells.Item(122, 14).Value = -13132   # N122: was -11382.4999
$ws.Cells.Item(136, 8).Value = 8184.9   # H136: was 7112.4165
$ws.Cells.Item(136, 9).Value = 4412.5   # I136: was 3810
$ws.Cells.Item(136, 10).Value = 10699.833   # J136: was 9471.286
$ws.Cells.Item(136, 11).Value = 13237.5   # K136: was 11430
$ws.Cells.Item(136, 12).Value = 32099.499   # L136: was 28413.858
$ws.Cells.Item(136, 13).Value = -10687.5   # M136: was -8880
$ws.Cells.Item(136, 14).Value = -37199.499   # N136: was -33513.858
# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 567.94116   # H2: was 517
$ws.Cells.Item(2, 9).Value = 598.125   # I2: was 567.8823
$ws.Cells.Item(2, 10).Value = 85   # J2: was 84.5
$ws.Cells.Item(2, 11).Value = 3588.75   # K2: was 3407.2938
$ws.Cells.Item(2, 12).Value = 510   # L2: was 507
$ws.Cells.Item(2, 13).Value = -3475.75   # M2: was -3294.2938
$ws.Cells.Item(2, 14).Value = -736   # N2: was -733
$ws.Cells.Item(5, 8).Value = 1346.95   # H5: was 1347
$ws.Cells.Item(5, 9).Value = 1164.85   # I5: was 1133.1904
$ws.Cells.Item(5, 10).Value = 1529.05   # J5: was 1583.3158
$ws.Cells.Item(5, 11).Value = 3494.55   # K5: was 3399.5712
$ws.Cells.Item(5, 12).Value = 4587.15   # L5: was 4749.9474
$ws.Cells.Item(5, 13).Value = -3382.55   # M5: was -3287.5712
$ws.Cells.Item(5, 14).Value = -4811.15   # N5: was -4973.9474
$ws.Cells.Item(58, 8).Value = 5000   # H58: was 0
$ws.Cells.Item(58, 10).Value = 5000   # J58: was 0
$ws.Cells.Item(58, 12).Value = 15000   # L58: was 0
$ws.Cells.Item(58, 14).Value = -15256   # N58: was None
$ws.Cells.Item(135, 8).Value = 1346.95   # H135: was 1347
$ws.Cells.Item(135, 9).Value = 1164.85   # I135: was 1133.1904
$ws.Cells.Item(135, 10).Value = 1529.05   # J135: was 1583.3158
$ws.Cells.Item(135, 11).Value = 10483.65   # K135: was 10198.7136
$ws.Cells.Item(135, 12).Value = 13761.45   # L135: was 14249.8422
$ws.Cells.Item(135, 13).Value = -7948.65   # M135: was -7663.713599999999
$ws.Cells.Item(135, 14).Value = -18831.45   # N135: was -19319.8422
# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 12682.45   # H70: was 12698.45
$ws.Cells.Item(70, 9).Value = 9630.888999999999   # I70: was 9666.444
$ws.Cells.Item(70, 11).Value = 9630.888999999999   # K70: was 9666.444
$ws.Cells.Item(70, 13).Value = -9360.888999999999   # M70: was -9396.444
$ws.Cells.Item(73, 8).Value = 12682.45   # H73: was 12698.45
$ws.Cells.Item(73, 9).Value = 9630.888999999999   # I73: was 9666.444
$ws.Cells.Item(73, 11).Value = 9630.888999999999   # K73: was 9666.444
$ws.Cells.Item(73, 13).Value = -8694.888999999999   # M73: was -8730.444
# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1283.875   # H22: was 1224.5555
$ws.Cells.Item(22, 9).Value = 1759   # I22: was 1506.75
$ws.Cells.Item(22, 11).Value = 1759   # K22: was 1506.75
$ws.Cells.Item(22, 13).Value = -1464   # M22: was -1211.75
$ws.Cells.Item(27, 8).Value = 1283.875   # H27: was 1224.5555
$ws.Cells.Item(27, 9).Value = 1759   # I27: was 1506.75
$ws.Cells.Item(27, 11).Value = 1759   # K27: was 1506.75
$ws.Cells.Item(27, 13).Value = -1652   # M27: was -1399.75
$ws.Cells.Item(46, 8).Value = 3299.6365   # H46: was 3108
$ws.Cells.Item(46, 9).Value = 1849.5   # I46: was 1679.6
$ws.Cells.Item(46, 11).Value = 1849.5   # K46: was 1679.6
$ws.Cells.Item(46, 13).Value = -1661.5   # M46: was -1491.6
$ws.Cells.Item(82, 8).Value = 11500.333   # H82: was 3566.8572
$ws.Cells.Item(82, 9).Value = 15000.5   # I82: was 3411.3333
$ws.Cells.Item(82, 11).Value = 15000.5   # K82: was 3411.3333
$ws.Cells.Item(82, 13).Value = -14639.5   # M82: was -3050.3333
$ws.Cells.Item(85, 8).Value = 11500.333   # H85: was 3566.8572
$ws.Cells.Item(85, 9).Value = 15000.5   # I85: was 3411.3333
$ws.Cells.Item(85, 11).Value = 15000.5   # K85: was 3411.3333
$ws.Cells.Item(85, 13).Value = -13752.5   # M85: was -2163.3333
$ws.Cells.Item(105, 8).Value = 7000   # H105: was 0
$ws.Cells.Item(105, 11).Value = 7000   # K105: was 0
$ws.Cells.Item(105, 12).Value = 7000   # L105: was 0
$ws.Cells.Item(105, 14).Value = -13988   # N105: was None
# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 6399   # H81: was 5666
$ws.Cells.Item(81, 9).Value = 2998.3333   # I81: was 2749
$ws.Cells.Item(81, 11).Value = 5996.6666   # K81: was 5498
$ws.Cells.Item(81, 13).Value = -4935.6666   # M81: was -4437
$ws.Cells.Item(84, 8).Value = 6399   # H84: was 5666
$ws.Cells.Item(84, 9).Value = 2998.3333   # I84: was 2749
$ws.Cells.Item(84, 11).Value = 29983.333   # K84: was 27490
$ws.Cells.Item(84, 13).Value = -24679.333   # M84: was -22186
$ws.Cells.Item(100, 8).Value = 1184.3334   # H100: was 1153.5
$ws.Cells.Item(100, 9).Value = 1206.7273   # I100: was 1163.75
$ws.Cells.Item(100, 11).Value = 2413.4546   # K100: was 2327.5
$ws.Cells.Item(100, 13).Value = -1872.4546   # M100: was -1786.5
